# Apply the supplier list correction:
#  - B2 (email of "supplier1") changes from "s@s.com" to "s1@s.com"
#  - Selection moves to E5 after the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the supplier1 email address
$ws.Range("B2").Value = "s1@s.com"

# Move / reflect the active selection to E5
$ws.Range("E5").Select()
